$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 384 (shifts rows 384..458 down to 385..459)
$ws.Rows.Item(384).Insert()

# Populate the newly inserted row 384 with the new data record
$ws.Cells.Item(384, 1).Value = 5
$ws.Cells.Item(384, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(384, 3).Value = "Maule"
$ws.Cells.Item(384, 4).Value = 44641
$ws.Cells.Item(384, 5).Value = 7
$ws.Cells.Item(384, 6).Value = "Fruta"
$ws.Cells.Item(384, 7).Value = 100102
$ws.Cells.Item(384, 8).Value = "Cítricos"
$ws.Cells.Item(384, 9).Value = 100102005
$ws.Cells.Item(384, 10).Value = "Naranja"
$ws.Cells.Item(384, 11).Value = "Valencia"
$ws.Cells.Item(384, 12).Value = "Primera"
$ws.Cells.Item(384, 13).Value = 360
$ws.Cells.Item(384, 14).Value = 10000
$ws.Cells.Item(384, 15).Value = 10000
$ws.Cells.Item(384, 16).Value = 10000
$ws.Cells.Item(384, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(384, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(384, 19).Value = 667
$ws.Cells.Item(384, 20).Value = 15
